$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
Write-Output $ws1.StandardWidth
